$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 647.2857
$ws.Range("I28").Value = 705.4706
$ws.Range("J28").Value = 400
$ws.Range("K28").Value = 705.4706
$ws.Range("L28").Value = 400
$ws.Range("M28").Value = -220.4706
$ws.Range("N28").Value = -1370
$ws.Range("H80").Value = 17857802
$ws.Range("I80").Value = 682.3333
$ws.Range("J80").Value = 31250640
$ws.Range("K80").Value = 2046.9999
$ws.Range("L80").Value = 93751920
$ws.Range("M80").Value = -1048.9999
$ws.Range("N80").Value = -93753916
$ws.Range("H83").Value = 17857802
$ws.Range("I83").Value = 682.3333
$ws.Range("J83").Value = 31250640
$ws.Range("K83").Value = 6140.9997
$ws.Range("L83").Value = 281255760
$ws.Range("M83").Value = -1148.9997
$ws.Range("N83").Value = -281265744
$ws.Range("H96").Value = 1230.25
$ws.Range("J96").Value = 1323.75
$ws.Range("L96").Value = 3971.25
$ws.Range("N96").Value = -6717.25
$ws.Range("H100").Value = 1327.6666
$ws.Range("I100").Value = 1327.6666
$ws.Range("K100").Value = 1327.6666
$ws.Range("M100").Value = -786.6666
$ws.Range("H123").Value = 109999.5
$ws.Range("J123").Value = 109999.5
$ws.Range("L123").Value = 109999.5
$ws.Range("N123").Value = -119799.5
$ws.Range("H135").Value = 363.03845
$ws.Range("I135").Value = 376.66666
$ws.Range("J135").Value = 199.5
$ws.Range("K135").Value = 3389.99994
$ws.Range("L135").Value = 1795.5
$ws.Range("M135").Value = -854.9999399999997
$ws.Range("N135").Value = -6865.5
$ws.Range("H137").Value = 3342558.5
$ws.Range("I137").Value = 6382.316
$ws.Range("J137").Value = 9105045
$ws.Range("K137").Value = 19146.948
$ws.Range("L137").Value = 27315135
$ws.Range("M137").Value = -16596.948
$ws.Range("N137").Value = -27320235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 252897.78
$ws.Range("I32").Value = 313028.2
$ws.Range("J32").Value = 12376.125
$ws.Range("K32").Value = 313028.2
$ws.Range("L32").Value = 12376.125
$ws.Range("M32").Value = -312741.2
$ws.Range("N32").Value = -12950.125
$ws.Range("H45").Value = 40572.54
$ws.Range("I45").Value = 45537.914
$ws.Range("K45").Value = 45537.914
$ws.Range("M45").Value = -45160.914
$ws.Range("H61").Value = 1357329.6
$ws.Range("I61").Value = 40965.465
$ws.Range("K61").Value = 40965.465
$ws.Range("M61").Value = -40753.465
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H136").Value = 1357329.6
$ws.Range("I136").Value = 40965.465
$ws.Range("K136").Value = 122896.395
$ws.Range("M136").Value = -120346.395

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 20756
$ws.Range("I107").Value = 20756
$ws.Range("K107").Value = 20756
$ws.Range("M107").Value = -18836
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3695.5
$ws.Range("I16").Value = 1397.5
$ws.Range("J16").Value = 5993.5
$ws.Range("K16").Value = 1397.5
$ws.Range("L16").Value = 5993.5
$ws.Range("M16").Value = -1110.5
$ws.Range("N16").Value = -6567.5
$ws.Range("H99").Value = 5001483
$ws.Range("I99").Value = 10001050
$ws.Range("K99").Value = 10001050
$ws.Range("M99").Value = -9999552
$ws.Range("H107").Value = 2258.5264
$ws.Range("J107").Value = 1897.5
$ws.Range("L107").Value = 1897.5
$ws.Range("N107").Value = -5737.5
$ws.Range("H113").Value = 3695.5
$ws.Range("I113").Value = 1397.5
$ws.Range("J113").Value = 5993.5
$ws.Range("K113").Value = 1397.5
$ws.Range("L113").Value = 5993.5
$ws.Range("M113").Value = 772.5
$ws.Range("N113").Value = -10333.5
$ws.Range("H126").Value = 5001483
$ws.Range("I126").Value = 10001050
$ws.Range("K126").Value = 30003150
$ws.Range("M126").Value = -30000680
$ws.Range("H132").Value = 3241.0833
$ws.Range("I132").Value = 4228
$ws.Range("K132").Value = 12684
$ws.Range("M132").Value = -10154
$ws.Range("H134").Value = 2803.1365
$ws.Range("I134").Value = 2952.6667
$ws.Range("J134").Value = 2482.7144
$ws.Range("K134").Value = 8858.000100000001
$ws.Range("L134").Value = 7448.1432
$ws.Range("M134").Value = -6323.000100000001
$ws.Range("N134").Value = -12518.1432
$ws.Range("H138").Value = 92917.69500000001
$ws.Range("J138").Value = 92917.69500000001
$ws.Range("L138").Value = 92917.69500000001
$ws.Range("N138").Value = -103197.695
$ws.Range("H139").Value = 61944.5
$ws.Range("J139").Value = 68000
$ws.Range("L139").Value = 68000
$ws.Range("N139").Value = -78280
$ws.Range("H141").Value = 264604
$ws.Range("J141").Value = 264604
$ws.Range("L141").Value = 264604
$ws.Range("N141").Value = -274964

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 208.55556
$ws.Range("I7").Value = 237.5
$ws.Range("K7").Value = 712.5
$ws.Range("M7").Value = -600.5
$ws.Range("H50").Value = 3828.182
$ws.Range("I50").Value = 3500.2
$ws.Range("J50").Value = 4101.5
$ws.Range("K50").Value = 10500.6
$ws.Range("L50").Value = 12304.5
$ws.Range("M50").Value = -10019.6
$ws.Range("N50").Value = -13266.5
$ws.Range("H53").Value = 3828.182
$ws.Range("I53").Value = 3500.2
$ws.Range("J53").Value = 4101.5
$ws.Range("K53").Value = 10500.6
$ws.Range("L53").Value = 12304.5
$ws.Range("M53").Value = -10019.6
$ws.Range("N53").Value = -13266.5
$ws.Range("H56").Value = 19236804
$ws.Range("I56").Value = 19236804
$ws.Range("K56").Value = 19236804
$ws.Range("M56").Value = -19236274
$ws.Range("H68").Value = 2602.9
$ws.Range("I68").Value = 3107
$ws.Range("K68").Value = 9321
$ws.Range("M68").Value = -8510
$ws.Range("H71").Value = 2602.9
$ws.Range("I71").Value = 3107
$ws.Range("K71").Value = 27963
$ws.Range("M71").Value = -23907
$ws.Range("H112").Value = 12008.667
$ws.Range("I112").Value = 10013.5
$ws.Range("K112").Value = 30040.5
$ws.Range("M112").Value = -28932.5
$ws.Range("H122").Value = 5556911
$ws.Range("I122").Value = 8333887.5
$ws.Range("J122").Value = 2958.25
$ws.Range("K122").Value = 75004987.5
$ws.Range("L122").Value = 26624.25
$ws.Range("M122").Value = -75002537.5
$ws.Range("N122").Value = -31524.25
$ws.Range("H131").Value = 4597921
$ws.Range("J131").Value = 74080.07000000001
$ws.Range("L131").Value = 222240.21
$ws.Range("N131").Value = -232320.21

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 27695.285
$ws.Range("J39").Value = 28811.166
$ws.Range("L39").Value = 28811.166
$ws.Range("N39").Value = -29875.166
$ws.Range("H80").Value = 25010352
$ws.Range("J80").Value = 45466870
$ws.Range("L80").Value = 45466870
$ws.Range("N80").Value = -45468866
$ws.Range("H83").Value = 25010352
$ws.Range("J83").Value = 45466870
$ws.Range("L83").Value = 227334350
$ws.Range("N83").Value = -227344334
$ws.Range("H102").Value = 20001534
$ws.Range("I102").Value = 23810852
$ws.Range("K102").Value = 23810852
$ws.Range("M102").Value = -23809230
$ws.Range("H132").Value = 1605878.2
$ws.Range("I132").Value = 1439.3334
$ws.Range("J132").Value = 2568541.5
$ws.Range("K132").Value = 4318.0002
$ws.Range("L132").Value = 7705624.5
$ws.Range("M132").Value = -1788.0002
$ws.Range("N132").Value = -7710684.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1589586.2
$ws.Range("I40").Value = 2224615
$ws.Range("K40").Value = 2224615
$ws.Range("M40").Value = -2224479
$ws.Range("H99").Value = 49508.25
$ws.Range("I99").Value = 36281.727
$ws.Range("J99").Value = 195000
$ws.Range("K99").Value = 36281.727
$ws.Range("L99").Value = 195000
$ws.Range("M99").Value = -33286.727
$ws.Range("N99").Value = -200990
$ws.Range("H136").Value = 55724.844
$ws.Range("I136").Value = 170416
$ws.Range("K136").Value = 511248
$ws.Range("M136").Value = -508698

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 545.3461
$ws.Range("I113").Value = 468.78946
$ws.Range("K113").Value = 1406.36838
$ws.Range("M113").Value = 763.6316199999999
$ws.Range("H114").Value = 68332.664
$ws.Range("J114").Value = 68332.664
$ws.Range("L114").Value = 68332.664
$ws.Range("N114").Value = -77010.664
$ws.Range("H132").Value = 2326.2068
$ws.Range("I132").Value = 2139
$ws.Range("J132").Value = 2526.7856
$ws.Range("K132").Value = 6417
$ws.Range("L132").Value = 7580.3568
$ws.Range("M132").Value = -3887
$ws.Range("N132").Value = -12640.3568
